# Insert the 2026 block of rows into Feuil1 (sheet1), mirroring the
# existing 2025 block (rows 198:206) one year later, with one content
# change in the last "rum" group (row 214): the regexp length moves
# from .{1,29} to .{1,32} and gets a distinguishing font/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the whole 2025 block (A198:E206) down to A207:E215 - this carries
# over all formatting (borders, number formats, shared-string reuse) so
# the new rows look exactly like their 2025 counterparts.
$srcBlock = $ws.Range("A198:E206")
$dstBlock = $ws.Range("A207:E215")
$srcBlock.Copy($dstBlock)

# Bump the year column for the whole new block to 2026.
$ws.Range("B207:B215").Value = 2026

# Row 214 (rum / zal group) differs from the 2025 pattern: the regexp
# length bucket changes from .{1,29} to .{1,32}.
$ws.Range("D214").Value = ".{1,32}"
$ws.Range("E214").Value = 32

# Give D214 a distinguishing font (new font + cell style), matching the
# extra font/cellXf introduced in this revision.
$ws.Range("D214").Font.Name = "Calibri"
$ws.Range("D214").Font.Color = 1

# Update the active selection to match where the edit was made.
$ws.Range("D213").Select() | Out-Null
